$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Cade Cunningham"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Detroit Pistons"

$ws.Range("A3").Value = "Derrick White"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Boston Celtics"

$ws.Range("A4").Value = "Jrue Holiday"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Boston Celtics"

$ws.Range("A5").Value = "Scotty Pippen Jr."
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Memphis Grizzlies"

$ws.Range("A6").Value = "Malik Monk"
$ws.Range("B6").Value = "PG,SG,SF"
$ws.Range("C6").Value = "Sacramento Kings"

$ws.Range("A7").Value = "Herbert Jones"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "New Orleans Pelicans"

$ws.Range("A8").Value = "Max Christie"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Los Angeles Lakers"

$ws.Range("A9").Value = "Isaiah Hartenstein"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Oklahoma City Thunder"

$ws.Range("A10").Value = "Anthony Davis"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Los Angeles Lakers"

$ws.Range("A11").Value = "Cason Wallace"
$ws.Range("B11").Value = "PG,SG"
$ws.Range("C11").Value = "Oklahoma City Thunder"

$ws.Range("A12").Value = "Damian Lillard"
$ws.Range("B12").Value = "PG"
$ws.Range("C12").Value = "Milwaukee Bucks"

$ws.Range("A13").Value = "Bam Adebayo"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Miami Heat"

$ws.Range("A14").Value = "Julius Randle"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Minnesota Timberwolves"

$ws.Range("A15").Value = "Donte DiVincenzo"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Minnesota Timberwolves"

$ws.Range("A16").Value = "Jaden McDaniels"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Minnesota Timberwolves"

$ws.Range("A17").Value = "Brandon Miller"
$ws.Range("B17").Value = "SG,SF,PF"
$ws.Range("C17").Value = "Charlotte Hornets"

$ws.Range("A18").Value = "LaMelo Ball"
$ws.Range("B18").Value = "PG,SG"
$ws.Range("C18").Value = "Charlotte Hornets"

$ws.Range("A19").Value = "Cameron Johnson"
$ws.Range("B19").Value = "SF,PF"
$ws.Range("C19").Value = "Brooklyn Nets"

